# Add three new rows (ARIMA, LSTM Univariate, ARIMA with seasonality) to the
# "Summary" results table that already lives on the first worksheet of the
# workbook (rows 2-9, columns B:K). The new rows continue the same layout and
# reuse the existing row-9 formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clone the formatting of the last existing data row (row 9) down into the
# three new rows so fills/borders/number formats line up with the rest of
# the table.
$ws.Range("B9:K9").Copy() | Out-Null
$ws.Range("B10:K12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 10: ARIMA ---------------------------------------------------
$ws.Range("B10").Value = "ARIMA"
# --- Row 11: LSTM Univariate ------------------------------------------
$ws.Range("B11").Value = "LSTM Univariate"
$ws.Range("G11").Value = "Test RMSE = 0.34"
$ws.Range("J11").Value = "Test RMSE = 0.33"
$ws.Range("D11").Value = "RMSE= 0.34"
# --- back to row 10 ----------------------------------------------------
$ws.Range("G10").Value = "RMSE =  0.21"
$ws.Range("J10").Value = "RMSE = 0.27"
# --- Row 12: ARIMA with seasonality ------------------------------------
$ws.Range("B12").Value = "ARIMA with seasonality"
$ws.Range("G12").Value = "RMSE =  0.39"

# Fill in the "N/A" placeholders for the Type columns and the one numeric
# score that stayed a number instead of becoming text.
$ws.Range("C10").Value = "N/A"
$ws.Range("F10").Value = "N/A"
$ws.Range("I10").Value = "N/A"
$ws.Range("D10").Value = 0.29
$ws.Range("F11").Value = "N/A"
$ws.Range("I11").Value = "N/A"
$ws.Range("C12").Value = "N/A"
$ws.Range("F12").Value = "N/A"
$ws.Range("I12").Value = "N/A"
$ws.Range("J12").Value = "RMSE = 0.27"

# J10 and J12 should carry the plain "N/A" style (general number format)
# rather than the percentage style copied from row 9, so pull that
# formatting over from the neighbouring "N/A" cell before writing the text.
$ws.Range("I10").Copy() | Out-Null
$ws.Range("J10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("J10").Value = "RMSE = 0.27"

$ws.Range("I12").Copy() | Out-Null
$ws.Range("J12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("J12").Value = "RMSE = 0.27"

# Match the author's final cursor position.
$ws.Range("D11").Select() | Out-Null
